# Apply updated cryptocurrency price/volume figures to sheet1 (cryptos.xlsx).
# Source data was regenerated by the scraping Action on Sat Sep  2 04:50:34 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.795.18"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.634.90"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'215.48"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'0.5029"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "'0.2576"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.06423"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'19.67"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "'0.07690"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.254"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.637.74"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "1.859.57"
$ws.Range("D15").Value = "'0.5461"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "0.0₅7943"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "25.826.11"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'203.77"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Value = "'4.332"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("D22").Value = "'9.958"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "'1.912"
$ws.Range("E25").Value = "  +8.99%  "
$ws.Range("D26").Value = "'141.28"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "'0.1147"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "'15.73"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "'0.05029"
$ws.Range("E30").Value = "  -3.45%  "
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "'3.273"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'3.187"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'1.535"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "1.176.23"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'0.8961"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("D39").Value = "'0.5612"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'0.01562"
$ws.Range("D41").Value = "'2.546"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "'5.674"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").Value = "'0.8081"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("D45").Value = "'99.64"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "1.771.89"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'0.4513"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'54.89"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").Value = "'0.05042"
$ws.Range("E51").Value = "  -0.41%  "
